$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal text into the cell even when it looks like a
    # number (e.g. "593.59", "1.00", "0.0000177"): briefly switch the
    # cell to Text format so the assignment is not reinterpreted as a
    # number, then clear the (temporary) formatting back off so the
    # cell's style is left exactly as it started (no explicit s=...).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '67.417.31'
$ws.Range("E2").Value = '  +0.74%  '
Set-TextValue $ws.Range("D3") '3.467.14'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue $ws.Range("D5") '593.59'
$ws.Range("E5").Value = '  +0.16%  '
Set-TextValue $ws.Range("D6") '181.44'
$ws.Range("E6").Value = '  +2.94%  '
$ws.Range("E7").Value = '  +5.21%  '
Set-TextValue $ws.Range("D8") '1.00'
$ws.Range("E8").Value = '  +0.15%  '
Set-TextValue $ws.Range("D9") '3.468.39'
$ws.Range("E9").Value = '  +0.24%  '
Set-TextValue $ws.Range("D10") '0.141'
$ws.Range("E10").Value = '  +9.19%  '
Set-TextValue $ws.Range("D11") '6.99'
$ws.Range("E11").Value = '  -1.12%  '
Set-TextValue $ws.Range("D12") '0.430'
$ws.Range("E12").Value = '  +0.87%  '
Set-TextValue $ws.Range("D13") '4.067.25'
$ws.Range("E13").Value = '  +0.42%  '
Set-TextValue $ws.Range("D14") '31.98'
$ws.Range("E14").Value = '  +3.99%  '
Set-TextValue $ws.Range("D15") '0.133'
$ws.Range("E15").Value = '  -0.62%  '
Set-TextValue $ws.Range("D16") '67.439.12'
$ws.Range("E16").Value = '  +0.84%  '
Set-TextValue $ws.Range("D17") '0.0000177'
$ws.Range("E17").Value = '  +1.08%  '
Set-TextValue $ws.Range("D18") '3.471.78'
$ws.Range("E18").Value = '  +1.20%  '
Set-TextValue $ws.Range("D19") '6.19'
$ws.Range("E19").Value = '  -0.51%  '
Set-TextValue $ws.Range("D20") '14.11'
$ws.Range("E20").Value = '  -1.25%  '
Set-TextValue $ws.Range("D21") '394.61'
$ws.Range("E21").Value = '  +1.66%  '
Set-TextValue $ws.Range("D22") '7.93'
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D24") '0.998'
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D25") '0.539'
$ws.Range("E25").Value = '  +1.17%  '
Set-TextValue $ws.Range("D26") '71.82'
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("E27").Value = '  +1.03%  '
Set-TextValue $ws.Range("D28") '10.34'
$ws.Range("E28").Value = '  +1.42%  '
Set-TextValue $ws.Range("D29") '0.175'
$ws.Range("E29").Value = '  -1.23%  '
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.15%  '
Set-TextValue $ws.Range("D31") '6.12'
$ws.Range("E31").Value = '  +0.64%  '
Set-TextValue $ws.Range("D32") '1.40'
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("E33").Value = '  +1.33%  '
Set-TextValue $ws.Range("D34") '23.53'
$ws.Range("E34").Value = '  +0.75%  '
Set-TextValue $ws.Range("D35") '7.33'
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("E36").Value = '  -0.16%  '
Set-TextValue $ws.Range("D37") '1.58'
$ws.Range("E37").Value = '  -1.98%  '
Set-TextValue $ws.Range("D38") '161.01'
$ws.Range("E38").Value = '  -1.10%  '
Set-TextValue $ws.Range("D39") '0.889'
$ws.Range("E39").Value = '  +2.70%  '
Set-TextValue $ws.Range("D40") '2.84'
$ws.Range("E40").Value = '  +11.64%  '
$ws.Range("E41").Value = '  -2.19%  '
Set-TextValue $ws.Range("D42") '6.77'
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("E43").Value = '  +1.46%  '
Set-TextValue $ws.Range("D44") '26.22'
$ws.Range("E44").Value = '  +0.51%  '
Set-TextValue $ws.Range("D45") '0.0717'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D46") '2.742.89'
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D47") '26.21'
$ws.Range("E47").Value = '  -3.41%  '
Set-TextValue $ws.Range("D48") '41.55'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("E49").Value = '  +0.23%  '
Set-TextValue $ws.Range("D50") '325.54'
$ws.Range("E50").Value = '  -4.21%  '
Set-TextValue $ws.Range("D51") '1.04'
$ws.Range("E51").Value = '  -2.04%  '
